# Sistema Drugs and Drage ons - "Alteração do diagrama de classe e inclusão do SOLID"
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 13 "Open-closed principle" - fix typo migraçÕes -> migrações
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$body13 = $s13.Shapes.Item(2)
$tr13 = $body13.TextFrame.TextRange
$full13 = $tr13.Text
$idx13 = $full13.IndexOf('migra')
$word13 = $tr13.Characters($idx13 + 1, 9)
$word13.Text = 'migrações'

# ---------------------------------------------------------------------------
# Helper: fill an (empty) TextRange with a sequence of text segments, each
# becoming its own <a:r> run, formatted as pt-BR, without leaving a stray
# endParaRPr behind.
# ---------------------------------------------------------------------------
function Set-MultiRunText {
    param(
        $TextRange,
        [string[]]$Parts
    )

    $joined = [string]::Join('', $Parts)
    $TextRange.Text = $joined
    $TextRange.LanguageID = 'pt-BR'

    $pos = 1
    foreach ($part in $Parts) {
        $len = $part.Length
        if ($len -gt 0) {
            $sub = $TextRange.Characters($pos, $len)
            $sub.Text = $part
        }
        $pos += $len
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 14 "Liskov substitution principle" - fill in the empty body
#    placeholder with the LSP paragraph.
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$body14 = $s14.Shapes.Item(2)
$tr14 = $body14.TextFrame.TextRange

$parts14 = @(
    'Em nosso trabalho, o LSP foi seguido e isso pode ser observado nas classes “Usuário”, “Cliente” e “Funcionário”. Anteriormente, a classe Usuário possuía um atributo “',
    'IDcliente',
    '” e "',
    'IDfuncionario',
    '", o que feria o Princípio de ',
    'Liskov',
    ', visto que não se podia substituir a classe “Usuário” por qualquer um de suas subclasses. Então, seguindo o princípio SOLID, ambos os atributos foram substituídos pelo atributo novo “',
    'IDusuário',
    '”.'
)

Set-MultiRunText $tr14 $parts14

# ---------------------------------------------------------------------------
# 3) Slide 15 "Interface segregation principle"
#    a) move/resize the title placeholder (explicit xfrm added)
#    b) fill in the empty body placeholder with the segregation paragraph
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$title15 = $s15.Shapes.Item(1)
$title15.Left = 24.543317086614174
$title15.Top = 35.62670291338583
$title15.Width = 670.9133958267716
$title15.Height = 55.7007974015748

$body15 = $s15.Shapes.Item(2)
$tr15 = $body15.TextFrame.TextRange

$parts15 = @(
    'Para o princípio da ',
    'segração',
    ' da interface, poderia ter sido utilizado uma interface de pagamento que permitisse a comunicação com qualquer sistema de validação de compras, sem prejuízo da completude das transações. Assim, foi criada uma interface de pagamento para realizar esta comunicação. De acordo com o padrão ',
    'Adapter',
    '.'
)

Set-MultiRunText $tr15 $parts15

Write-Output 'done'
